$d = $word.ActiveDocument

$sentence = "Now the sub goals are the delivery of each parcel individually while leaving two parcels unprotected at two given times."
$newSentence = "So how do we do this? If you transport the seed this would leave the cat and parrot alone for destruction. If you transport the cat this leaves the seed in jeopardy with the parrot. When you transport the parrot the seed and cat are safe. But when return to transport and deliver the either the seed or cat we have problems on the other side of the river as they are left alone without the protector. I am assuming that the man needs to be presented on this boat at all times. Otherwise who would guide this boat from point A to point B?"

# Locate the (currently unique) sentence "Now the sub goals..." - its
# paragraph currently reads: <tab> + <bookmark> + <sentence>.
$r = $d.Content
$r.Find.Execute($sentence, $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r.Collapse(1)
$tabParaIndex = $r.Paragraphs(1).Index

# Split right before the sentence: this leaves the leading tab alone in
# its own paragraph ($tabParaIndex), and pushes the bookmark + sentence
# into the following paragraph ($tabParaIndex + 1).
$r.InsertParagraphBefore()

# Re-create "Now the sub goals..." after the tab, ahead of the bookmark.
$tabPara = $d.Paragraphs($tabParaIndex)
$tabPara.Range.InsertAfter($sentence)

# The bookmark (and the now-duplicate original sentence run) live in the
# next paragraph. Insert a blank paragraph, then a new paragraph (tab +
# the new sentence), both ahead of the bookmark's paragraph.
$bookmarkParaIndex = $tabParaIndex + 1

$bookmarkPara = $d.Paragraphs($bookmarkParaIndex)
$insertPoint = $bookmarkPara.Range.Duplicate
$insertPoint.Collapse(1)
$insertPoint.InsertParagraphBefore()
# Blank paragraph is now at $bookmarkParaIndex; bookmark paragraph moved
# to $bookmarkParaIndex + 1.

$bookmarkPara = $d.Paragraphs($bookmarkParaIndex + 1)
$insertPoint2 = $bookmarkPara.Range.Duplicate
$insertPoint2.Collapse(1)
$insertPoint2.InsertParagraphBefore()
# New (currently empty) paragraph is now at $bookmarkParaIndex + 1;
# bookmark paragraph moved to $bookmarkParaIndex + 2.

$newTextPara = $d.Paragraphs($bookmarkParaIndex + 1)
$newTextPara.Range.InsertAfter("`t" + $newSentence)

# Finally, remove the leftover duplicate sentence text sitting beside the
# bookmark, leaving the bookmark alone in its paragraph.
$bookmarkPara = $d.Paragraphs($bookmarkParaIndex + 2)
$bookmarkPara.Range.Text = ""
